$d = $word.ActiveDocument
$d.Content.Find.Execute("SPO 5", $true, $false, $false, $false, $false, $true, 1, $false, "SPO 4", 2)
